# Correlation4Versioning.xlsx edit:
#  - fix the "WSGC02 WSGC02" typo in C2 -> "WSGC01 WSGC02 WSGC03"
#  - re-sort the commentary list in column B (rows 27-33) so the two
#    1888-08-26 / 1888-08-27 dates sit at the top of that block
#  - drop the two trailing blank placeholder rows (34-35) and recreate
#    them two rows further down (36-37) now that the used range grew
#  - leave the selection on C2, matching the saved view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal/text value to a cell without Excel's COM layer
# auto-coercing date-shaped strings (e.g. "1888-08-26") into date serials.
# A leading apostrophe forces text entry, exactly like typing it in the UI.
function Set-TextValue($rng, [string]$val) {
    $rng.Value = "'" + $val
}

# --- 1. Typo fix: WSGC01 WSGC02 WSGC02 -> WSGC01 WSGC02 WSGC03 ---------
$ws.Range("C2").Value = "WSGC01 WSGC02 WSGC03"

# --- 2. Remove the two blank, style-only rows at the bottom (34:35); ---
#        they get re-added (shifted down by 2) once the B-column list
#        below grows by two rows.
[void]$ws.Rows("34:35").Delete()

# --- 3. Re-order the B27:B33 block -------------------------------------
Set-TextValue $ws.Range("B27") "1888-08-26"
Set-TextValue $ws.Range("B28") "1888-08-27"
Set-TextValue $ws.Range("B29") "Milton George and Charles H. Ham Commentary:1888-08-28"
Set-TextValue $ws.Range("B30") "Judge O.H. Horton Commentary:1888-08-30"
Set-TextValue $ws.Range("B31") "Col. Abner Taylor Commentary:1888-08-31"
Set-TextValue $ws.Range("B32") "Lyman J. Cage Commentary:1888-09-01"
Set-TextValue $ws.Range("B33") "George M. Sloan Commentary:1888-09-03"

# Re-apply the original column-B cell format (font/alignment/number format)
# that the apostrophe-prefix trick stamps with a stray quote-prefix flag;
# B4 already carries the same plain "General" text style every other
# B-column data cell uses.
[void]$ws.Range("B4").Copy()
[void]$ws.Range("B27:B33").PasteSpecial(-4122)

# --- 4. Recreate the two trailing placeholder rows at 36:37 ------------
[void]$ws.Range("A33").Copy()
[void]$ws.Range("A36:A37").PasteSpecial(-4122)

# --- 5. Selection / view moves to C2 ------------------------------------
[void]$ws.Range("C2").Select()

$excel.CutCopyMode = $false
